$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..11 (Colaborador_id, Colaborador_nome, Departamento, Motivo_da_ausência, Horas_de_ausência, Data_da_ausência, Salário)
$data = @(
    @(55965, "Dr. Breno Mendes",     "P&D",                     "Doença",             2, 45099, 4634.34),
    @(56263, "Mirella da Luz",       "Atendimento ao Cliente",  "Problemas pessoais", 5, 45096, 3525.39),
    @(43094, "Dr. Thomas das Neves", "Engenharia",               "Outros",             1, 45085, 4828.67),
    @(18717, "Benício Ribeiro",      "Atendimento ao Cliente",  "Doença",             4, 45087, 11598.63),
    @(75383, "Laura Azevedo",        "Operações",                "Consulta médica",    6, 45099, 4089.4),
    @(92050, "Melissa da Costa",     "TI",                       "Viagem de negócios", 8, 45096, 5015.86),
    @(86348, "Breno Duarte",         "TI",                       "Problemas pessoais", 5, 45080, 3857.9),
    @(32577, "Dra. Sophie Duarte",   "Operações",                "Problemas pessoais", 4, 45090, 8672.19),
    @(30843, "Felipe Almeida",       "Atendimento ao Cliente",  "Outros",             1, 45092, 7355.54),
    @(68333, "Sarah da Rocha",       "Engenharia",               "Viagem de negócios", 6, 45090, 11117.24)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
